# Apply the edits described by the diff:
# 1. Fix the shared string "MODEL_CONDITION" -> "MODELCONDITION"
#    (this is just the text of whatever cell held that value; it will be
#    found in what becomes column D after the column deletion below)
# 2. Delete column A entirely, shifting B:F left to A:E

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A (the taxon-count helper column), shifting remaining
# columns (B:F) left to (A:E).
$ws.Range("A:A").Delete()

# Correct the header text that used to read "MODEL_CONDITION".
# After the column shift this header now lives in column D.
$ws.Range("D1").Value = "MODELCONDITION"
